$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7692.3076
$ws.Range("I69").Value = 14000
$ws.Range("J69").Value = 4888.8887
$ws.Range("K69").Value = 42000
$ws.Range("L69").Value = 14666.6661
$ws.Range("M69").Value = -41126
$ws.Range("N69").Value = -16414.6661

$ws.Range("H72").Value = 7692.3076
$ws.Range("I72").Value = 14000
$ws.Range("J72").Value = 4888.8887
$ws.Range("K72").Value = 126000
$ws.Range("L72").Value = 43999.99830000001
$ws.Range("M72").Value = -121632
$ws.Range("N72").Value = -52735.99830000001

$ws.Range("H74").Value = 4674.875
$ws.Range("I74").Value = 4674.875
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4674.875
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3738.875
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 4674.875
$ws.Range("I77").Value = 4674.875
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 23374.375
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -18694.375
$ws.Range("N77").Value = -18694.375

$ws.Range("H113").Value = 5557404.5
$ws.Range("I113").Value = 8335025.5
$ws.Range("K113").Value = 8335025.5
$ws.Range("M113").Value = -8331771.5

$ws.Range("H116").Value = 13740745
$ws.Range("I116").Value = 8335999.5
$ws.Range("J116").Value = 18544964
$ws.Range("K116").Value = 8335999.5
$ws.Range("L116").Value = 18544964
$ws.Range("M116").Value = -8332557.5
$ws.Range("N116").Value = -18551848

$ws.Range("H129").Value = 479217.56
$ws.Range("I129").Value = 1605.7778
$ws.Range("J129").Value = 837426.4399999999
$ws.Range("K129").Value = 4817.3334
$ws.Range("L129").Value = 2512279.32
$ws.Range("M129").Value = 182.6665999999996
$ws.Range("N129").Value = -2522279.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 25312.268
$ws.Range("I2").Value = 30101.75
$ws.Range("J2").Value = 6154.3335
$ws.Range("K2").Value = 30101.75
$ws.Range("L2").Value = 6154.3335
$ws.Range("M2").Value = -29988.75
$ws.Range("N2").Value = -6380.3335

$ws.Range("H32").Value = 2163106.2
$ws.Range("I32").Value = 2785346
$ws.Range("J32").Value = 9199.23
$ws.Range("K32").Value = 2785346
$ws.Range("L32").Value = 9199.23
$ws.Range("M32").Value = -2785059
$ws.Range("N32").Value = -9773.23

$ws.Range("H97").Value = 473.3889
$ws.Range("I97").Value = 467.33334
$ws.Range("J97").Value = 503.66666
$ws.Range("K97").Value = 467.33334
$ws.Range("L97").Value = 503.66666
$ws.Range("M97").Value = 28.66665999999998
$ws.Range("N97").Value = -1495.66666

$ws.Range("H110").Value = 759.625
$ws.Range("I110").Value = 698.75
$ws.Range("K110").Value = 698.75
$ws.Range("M110").Value = 1346.25

$ws.Range("H116").Value = 25312.268
$ws.Range("I116").Value = 30101.75
$ws.Range("J116").Value = 6154.3335
$ws.Range("K116").Value = 30101.75
$ws.Range("L116").Value = 6154.3335
$ws.Range("M116").Value = -27807.75
$ws.Range("N116").Value = -10742.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 25312.268
$ws.Range("I3").Value = 30101.75
$ws.Range("J3").Value = 6154.3335
$ws.Range("K3").Value = 30101.75
$ws.Range("L3").Value = 6154.3335
$ws.Range("M3").Value = -29987.75
$ws.Range("N3").Value = -6382.3335

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H86").Value = 1963.92
$ws.Range("I86").Value = 1963.5858
$ws.Range("K86").Value = 1963.5858
$ws.Range("M86").Value = -840.5858000000001

$ws.Range("H89").Value = 1963.92
$ws.Range("I89").Value = 1963.5858
$ws.Range("K89").Value = 9817.929
$ws.Range("M89").Value = -4201.929

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4226.9614
$ws.Range("I122").Value = 6288.5
$ws.Range("K122").Value = 18865.5
$ws.Range("M122").Value = -16415.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3122113.5
$ws.Range("I5").Value = 2137198.8
$ws.Range("J5").Value = 6667806
$ws.Range("K5").Value = 6411596.399999999
$ws.Range("L5").Value = 20003418
$ws.Range("M5").Value = -6411484.399999999
$ws.Range("N5").Value = -20003642

$ws.Range("H131").Value = 53492.473
$ws.Range("J131").Value = 126374.5
$ws.Range("L131").Value = 379123.5
$ws.Range("N131").Value = -389203.5

$ws.Range("H135").Value = 3122113.5
$ws.Range("I135").Value = 2137198.8
$ws.Range("J135").Value = 6667806
$ws.Range("K135").Value = 19234789.2
$ws.Range("L135").Value = 60010254
$ws.Range("M135").Value = -19232254.2
$ws.Range("N135").Value = -60015324

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4404282
$ws.Range("I122").Value = 31060
$ws.Range("J122").Value = 9806498
$ws.Range("K122").Value = 93180
$ws.Range("L122").Value = 29419494
$ws.Range("M122").Value = -90730
$ws.Range("N122").Value = -29424394

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2086.1538
$ws.Range("I68").Value = 1787.6923
$ws.Range("J68").Value = 2384.6155
$ws.Range("K68").Value = 1787.6923
$ws.Range("L68").Value = 2384.6155
$ws.Range("M68").Value = -1038.6923
$ws.Range("N68").Value = -3882.6155

$ws.Range("H71").Value = 2086.1538
$ws.Range("I71").Value = 1787.6923
$ws.Range("J71").Value = 2384.6155
$ws.Range("K71").Value = 8938.461499999999
$ws.Range("L71").Value = 11923.0775
$ws.Range("M71").Value = -5194.461499999999
$ws.Range("N71").Value = -19411.0775

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22740786
$ws.Range("I62").Value = 55580876
$ws.Range("J62").Value = 5336.846
$ws.Range("K62").Value = 55580876
$ws.Range("L62").Value = 5336.846
$ws.Range("M62").Value = -55580252
$ws.Range("N62").Value = -6584.846

$ws.Range("H65").Value = 22740786
$ws.Range("I65").Value = 55580876
$ws.Range("J65").Value = 5336.846
$ws.Range("K65").Value = 277904380
$ws.Range("L65").Value = 26684.23
$ws.Range("M65").Value = -277901260
$ws.Range("N65").Value = -32924.23

$ws.Range("H81").Value = 29610.4
$ws.Range("I81").Value = 1300
$ws.Range("J81").Value = 32756
$ws.Range("K81").Value = 2600
$ws.Range("L81").Value = 65512
$ws.Range("M81").Value = -1539
$ws.Range("N81").Value = -67634

$ws.Range("H84").Value = 29610.4
$ws.Range("I84").Value = 1300
$ws.Range("J84").Value = 32756
$ws.Range("K84").Value = 13000
$ws.Range("L84").Value = 327560
$ws.Range("M84").Value = -7696
$ws.Range("N84").Value = -338168

$ws.Range("H107").Value = 18191.727
$ws.Range("I107").Value = 28272.715
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 84818.145
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = -82898.145
$ws.Range("N107").Value = -5490
